$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Col3a1"
$ws.Range("C2").Value = "Ddr2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.822558333333333
$ws.Range("H2").Value = 23.467675
$ws.Range("I2").Value = 0.001247993910151231
$ws.Range("J2").Value = 0.001247993910151231
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.4763
$ws.Range("N2").Value = 4.428900000000001
$ws.Range("O2").Value = 0.01318769285519422
$ws.Range("P2").Value = 0.01318769285519422
$ws.Range("Q2").Value = 11.5484428675
$ws.Range("R2").Value = 103.9359858075
$ws.Range("S2").Value = 0.00001645816037222729
$ws.Range("T2").Value = 0.00001645816037222729

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Col3a1"
$ws.Range("C3").Value = "Ddr2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.822558333333333
$ws.Range("H3").Value = 23.467675
$ws.Range("I3").Value = 0.001247993910151231
$ws.Range("J3").Value = 0.001247993910151231
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 78.17189533333334
$ws.Range("N3").Value = 234.515686
$ws.Range("O3").Value = 0.6983045082736506
$ws.Range("P3").Value = 0.6983045082736506
$ws.Range("Q3").Value = 611.5042112722278
$ws.Range("R3").Value = 5503.53790145005
$ws.Range("S3").Value = 0.0008714797737566659
$ws.Range("T3").Value = 0.0008714797737566659

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Col3a1"
$ws.Range("C4").Value = "Ddr2"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.822558333333333
$ws.Range("H4").Value = 23.467675
$ws.Range("I4").Value = 0.001247993910151231
$ws.Range("J4").Value = 0.001247993910151231
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1536526666666667
$ws.Range("N4").Value = 0.460958
$ws.Range("O4").Value = 0.00137256937911098
$ws.Range("P4").Value = 0.00137256937911098
$ws.Range("Q4").Value = 1.201956948072222
$ws.Range("R4").Value = 10.81761253265
$ws.Range("S4").Value = 0.000001712958226390559
$ws.Range("T4").Value = 0.000001712958226390559

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Col3a1"
$ws.Range("C5").Value = "Ddr2"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 7.822558333333333
$ws.Range("H5").Value = 23.467675
$ws.Range("I5").Value = 0.001247993910151231
$ws.Range("J5").Value = 0.001247993910151231
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 32.14343433333334
$ws.Range("N5").Value = 96.43030300000001
$ws.Range("O5").Value = 0.2871352294920441
$ws.Range("P5").Value = 0.2871352294920441
$ws.Range("Q5").Value = 251.4438901061695
$ws.Range("R5").Value = 2262.995010955525
$ws.Range("S5").Value = 0.0003583430177959471
$ws.Range("T5").Value = 0.0003583430177959471

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Col3a1"
$ws.Range("C6").Value = "Ddr2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 6037.773437333333
$ws.Range("H6").Value = 18113.320312
$ws.Range("I6").Value = 0.9632532171165058
$ws.Range("J6").Value = 0.9632532171165058
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.4763
$ws.Range("N6").Value = 4.428900000000001
$ws.Range("O6").Value = 0.01318769285519422
$ws.Range("P6").Value = 0.01318769285519422
$ws.Range("Q6").Value = 8913.564925535202
$ws.Range("R6").Value = 80222.0843298168
$ws.Range("S6").Value = 0.01270308756911019
$ws.Range("T6").Value = 0.01270308756911019

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Col3a1"
$ws.Range("C7").Value = "Ddr2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 6037.773437333333
$ws.Range("H7").Value = 18113.320312
$ws.Range("I7").Value = 0.9632532171165058
$ws.Range("J7").Value = 0.9632532171165058
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 78.17189533333334
$ws.Range("N7").Value = 234.515686
$ws.Range("O7").Value = 0.6983045082736506
$ws.Range("P7").Value = 0.6983045082736506
$ws.Range("Q7").Value = 471984.1931896016
$ws.Range("R7").Value = 4247857.738706415
$ws.Range("S7").Value = 0.6726440641215536
$ws.Range("T7").Value = 0.6726440641215536

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Col3a1"
$ws.Range("C8").Value = "Ddr2"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 6037.773437333333
$ws.Range("H8").Value = 18113.320312
$ws.Range("I8").Value = 0.9632532171165058
$ws.Range("J8").Value = 0.9632532171165058
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1536526666666667
$ws.Range("N8").Value = 0.460958
$ws.Range("O8").Value = 0.00137256937911098
$ws.Range("P8").Value = 0.00137256937911098
$ws.Range("Q8").Value = 927.7199893754329
$ws.Range("R8").Value = 8349.479904378895
$ws.Range("S8").Value = 0.001322131870144256
$ws.Range("T8").Value = 0.001322131870144256

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Col3a1"
$ws.Range("C9").Value = "Ddr2"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 6037.773437333333
$ws.Range("H9").Value = 18113.320312
$ws.Range("I9").Value = 0.9632532171165058
$ws.Range("J9").Value = 0.9632532171165058
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 32.14343433333334
$ws.Range("N9").Value = 96.43030300000001
$ws.Range("O9").Value = 0.2871352294920441
$ws.Range("P9").Value = 0.2871352294920441
$ws.Range("Q9").Value = 194074.7740024683
$ws.Range("R9").Value = 1746672.966022215
$ws.Range("S9").Value = 0.2765839335556977
$ws.Range("T9").Value = 0.2765839335556977

# Row 10
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Col3a1"
$ws.Range("C10").Value = "Ddr2"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.018774
$ws.Range("H10").Value = 3.056322
$ws.Range("I10").Value = 0.0001625329839219791
$ws.Range("J10").Value = 0.0001625329839219791
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 1.4763
$ws.Range("N10").Value = 4.428900000000001
$ws.Range("O10").Value = 0.01318769285519422
$ws.Range("P10").Value = 0.01318769285519422
$ws.Range("Q10").Value = 1.5040160562
$ws.Range("R10").Value = 13.5361445058
$ws.Range("S10").Value = 0.000002143435070801281
$ws.Range("T10").Value = 0.000002143435070801281

# Row 11
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Col3a1"
$ws.Range("C11").Value = "Ddr2"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.018774
$ws.Range("H11").Value = 3.056322
$ws.Range("I11").Value = 0.0001625329839219791
$ws.Range("J11").Value = 0.0001625329839219791
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 78.17189533333334
$ws.Range("N11").Value = 234.515686
$ws.Range("O11").Value = 0.6983045082736506
$ws.Range("P11").Value = 0.6983045082736506
$ws.Range("Q11").Value = 79.63949449632133
$ws.Range("R11").Value = 716.755450466892
$ws.Range("S11").Value = 0.0001134975154158868
$ws.Range("T11").Value = 0.0001134975154158868

# Row 12
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Col3a1"
$ws.Range("C12").Value = "Ddr2"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.018774
$ws.Range("H12").Value = 3.056322
$ws.Range("I12").Value = 0.0001625329839219791
$ws.Range("J12").Value = 0.0001625329839219791
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1536526666666667
$ws.Range("N12").Value = 0.460958
$ws.Range("O12").Value = 0.00137256937911098
$ws.Range("P12").Value = 0.00137256937911098
$ws.Range("Q12").Value = 0.1565373418306666
$ws.Range("R12").Value = 1.408836076476
$ws.Range("S12").Value = 0.0000002230877968268457
$ws.Range("T12").Value = 0.0000002230877968268457

# Row 13
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Col3a1"
$ws.Range("C13").Value = "Ddr2"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.018774
$ws.Range("H13").Value = 3.056322
$ws.Range("I13").Value = 0.0001625329839219791
$ws.Range("J13").Value = 0.0001625329839219791
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 32.14343433333334
$ws.Range("N13").Value = 96.43030300000001
$ws.Range("O13").Value = 0.2871352294920441
$ws.Range("P13").Value = 0.2871352294920441
$ws.Range("Q13").Value = 32.74689516950733
$ws.Range("R13").Value = 294.722056525566
$ws.Range("S13").Value = 0.00004666894563846418
$ws.Range("T13").Value = 0.00004666894563846418

# Row 14
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Col3a1"
$ws.Range("C14").Value = "Ddr2"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 221.4914043333333
$ws.Range("H14").Value = 664.474213
$ws.Range("I14").Value = 0.03533625598942085
$ws.Range("J14").Value = 0.03533625598942085
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 1.4763
$ws.Range("N14").Value = 4.428900000000001
$ws.Range("O14").Value = 0.01318769285519422
$ws.Range("P14").Value = 0.01318769285519422
$ws.Range("Q14").Value = 326.9877602173
$ws.Range("R14").Value = 2942.8898419557
$ws.Range("S14").Value = 0.0004660036906409993
$ws.Range("T14").Value = 0.0004660036906409993

# Row 15
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Col3a1"
$ws.Range("C15").Value = "Ddr2"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 221.4914043333333
$ws.Range("H15").Value = 664.474213
$ws.Range("I15").Value = 0.03533625598942085
$ws.Range("J15").Value = 0.03533625598942085
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 78.17189533333334
$ws.Range("N15").Value = 234.515686
$ws.Range("O15").Value = 0.6983045082736506
$ws.Range("P15").Value = 0.6983045082736506
$ws.Range("Q15").Value = 17314.40287677835
$ws.Range("R15").Value = 155829.6258910051
$ws.Range("S15").Value = 0.02467546686292437
$ws.Range("T15").Value = 0.02467546686292437

# Row 16
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Col3a1"
$ws.Range("C16").Value = "Ddr2"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 221.4914043333333
$ws.Range("H16").Value = 664.474213
$ws.Range("I16").Value = 0.03533625598942085
$ws.Range("J16").Value = 0.03533625598942085
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1536526666666667
$ws.Range("N16").Value = 0.460958
$ws.Range("O16").Value = 0.00137256937911098
$ws.Range("P16").Value = 0.00137256937911098
$ws.Range("Q16").Value = 34.03274491956155
$ws.Range("R16").Value = 306.2947042760539
$ws.Range("S16").Value = 0.000048501462943506
$ws.Range("T16").Value = 0.000048501462943506

# Row 17
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Col3a1"
$ws.Range("C17").Value = "Ddr2"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 221.4914043333333
$ws.Range("H17").Value = 664.474213
$ws.Range("I17").Value = 0.03533625598942085
$ws.Range("J17").Value = 0.03533625598942085
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 32.14343433333334
$ws.Range("N17").Value = 96.43030300000001
$ws.Range("O17").Value = 0.2871352294920441
$ws.Range("P17").Value = 0.2871352294920441
$ws.Range("Q17").Value = 7119.494410586283
$ws.Range("R17").Value = 64075.44969527654
$ws.Range("S17").Value = 0.01014628397291197
$ws.Range("T17").Value = 0.01014628397291197
